$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 384 (the post "「いつの日か、私は飛ぶんだ」..."),
# which shifts all subsequent rows up by one.
$ws.Rows("384:384").Delete()
